$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $result = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $result) {
        Write-Output "FAILED to replace: $old"
    }
}

# Ativacao date 2018 -> 2022
Replace-Text "Ativação: 01/01/2018" "Ativação: 01/01/2022"

# Objetivos (PT)
Replace-Text "Apresentar ao aluno as características da profissão, os conceitos básicos envolvidos e as áreas de atuação do Engenheiro Bioquímico, e um panorama do mercado de trabalho de engenharia no Brasil. Mostrar como funciona uma indústria de bioprocesso, envolvendo instalações, tipos de bioprocessos e escala de produção." "Apresentar aos alunos a Engenharia Bioquímica, as características da profissão e orientar quanto as atribuições e as áreas de atuação do Engenheiro Bioquímico. Além disso, desenvolver nos alunos uma visão macro dos tipos e etapas de um bioprocesso industrial e, por fim, orientar sobre a atuação do Engenheiro Bioquímico na indústria, pesquisa e ensino, e empreendedorismo e inovação em engenharia."

# Objetivos (EN)
Replace-Text "Present to the student the characteristics of the profession, the basic concepts involved and the areas of expertise of the Biochemical Engineer, and an overview of the engineering work market in Brazil. Show how a bioprocessing industry works, involving facilities, types of bioprocesses and production scale." "To present to the Biochemical Engineering student the characteristics of the profession and to guide in relation to the attributes and the action areas of the biochemical engineering. Besides, to develop in the students a macro view of types and stages of an industrial bioprocess and, finally, to guide about the action of the biochemical engineering on the industry, research and teaching, and entrepreneurship and innovation in engineering."

# --- Add a second "Docente(s) Responsável(eis)" line: "5817181 - Valdeir Arantes" ---
# Locate the paragraph containing the first docente entry.
$docenteParaIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -like "101761*Arnaldo*Prata*") {
        $docenteParaIndex = $i
        break
    }
}
if ($docenteParaIndex -eq -1) {
    Write-Output "FAILED to find docente paragraph"
} else {
    $p = $d.Paragraphs.Item($docenteParaIndex)
    $r = $p.Range
    # Append a brand-new paragraph right after, holding the new docente text
    # (this keeps each run clean, with no leftover run-formatting markers).
    $insertPoint = $d.Range($r.End - 1, $r.End - 1)
    $insertPoint.InsertParagraphAfter()

    $full = $d.Content.Text
    $mark = $full.IndexOf("Prata") + 5

    $newParaStart = $d.Range($mark + 1, $mark + 1)
    $newParaStart.InsertAfter("5817181 - Valdeir Arantes")

    # Turn the trailing paragraph mark (right after "Prata") into a manual
    # line break, then delete the now-duplicate paragraph mark so the two
    # paragraphs merge back into one (two runs, split by <w:br/>).
    $breakPoint = $d.Range($mark, $mark)
    $breakPoint.Text = [string][char]11

    $mergeMark = $d.Range($mark + 1, $mark + 2)
    $mergeMark.Delete()
}

# Programa resumido / Programa text blocks are untouched by this diff,
# only the detailed "Programa" section below changes.

# Programa (PT)
Replace-Text "1. Histórico da Engenharia Bioquímica: interação entre ciências biológicas e a engenha, multidisciplinaridade, peculiaridades dos processos biotecnológicos. 2. Mercado de trabalho da Engenharia do Brasil3. Definições e conceitos – processo enzimático, processo fermentativo genérico, agentes de transformação, biorreator, matéria prima, tipos de substratos, conversão de substrato em produto, tipos de produtos biotecnológicos, recuperação de produtos, entre outros.4. Áreas de atuação do Engenheiro Bioquímico5. A Indústria de Bioprocessos – tipos de indústrias, equipamentos, instalações, principais operações unitárias. 6. Escalas de produção – laboratório, piloto, industrial. 7. Estudo de casos (processos biotecnológicos).8. Visitas supervisionadas – visitas a laboratórios e a indústria de bioprocesso." "1.Histórico da Engenharia Bioquímica: interação entre ciências biológicas e a engenha, multidisciplinaridade, peculiaridades dos processos biotecnológicos. 2.Mercado de trabalho da Engenharia do Brasil 3.Atribuições e áreas de atuação do Engenheiro Bioquímico 4.Definições e conceitos – processo enzimático, processo fermentativo genérico, agentes de transformação, biorreator, matéria prima, tipos de substratos, conversão de substrato em produto, tipos de produtos biotecnológicos, recuperação de produtos, entre outros. 5.A Indústria de Bioprocessos – tipos de indústrias, equipamentos, instalações, principais operações unitárias. 6.Escalas de produção – laboratório, piloto, industrial. 7.Estudo de casos (processos biotecnológicos). 8.Empreendedorismo e Inovação em Engenharia.9.Visitas supervisionadas – visitas a laboratórios e a indústria de bioprocesso."

# Programa (EN)
Replace-Text "1. History of Biochemical Engineering: interaction between biological sciences and engineering, multidisciplinarity, peculiarities of biotechnological processes.2. Labor market of Engineering in Brazil3. Definitions and concepts - enzymatic process, generic fermentation process, transformation agents, bioreactor, raw material, types of substrates, substrate conversion into product, types of biotechnological products, product recovery, among others.4. Areas of practice of the Biochemical Engineer5. The Bioprocess Industry - types of industries, equipment, facilities, main unit operations.6. Production scales - laboratory, pilot, industrial.7. Case studies (biotechnological processes).8. Supervised visits - visits to laboratories and the bioprocess industry." "1.History of the Biochemical Engineering: interaction between biological sciences and engineering, multidisciplinarity, peculiarities of biotechnological processes.2.Job market of Engineering in Brazil3.Attributes and action areas of biochemical engineering4.Definitions and concepts – enzymatic process, general fermentative process, transformation agents, bioreactor, raw material, types of substrates, conversion of substrate into product, types of biotechnological products, products recovery, between others.5.The Bioprocesses Industry – types of industries, equipment, installations, main unit operations6.Production scales – laboratory, pilot, industrial.7.Studies of cases (biotechnological processes).8.Entrepreneurship and Innovation in Engineering.9.Supervised visitation – visits to laboratories and bioprocess industry"

# Método
Replace-Text "Provas escritas; participação e conteúdo de trabalho e seminário;" "O método utilizado tem por fundamento a Aprendizagem Baseada em Projetos (PBL) que visa desenvolver as competências técnicas relativas ao tema do projeto, bem como competências transversais, tais como: aprender a aprender, trabalho em equipe, relacionamento interpessoal, aspectos de liderança e capacidade de comunicação, dentre outras; exercícios individuais realizados no decorrer da disciplina; exercícios; dinâmicas. Para os projetos, os alunos serão divididos em grupos que desenvolverão um projeto durante o semestre relacionado a aplicações dos conceitos abordados à um processo, produto ou serviço na área de Engenharia de Bioquímica e que relacione com a formação acadêmica e atribuições profissionais do Engenheiro Bioquímico."

# Critério
Replace-Text "A nota (N) será composta por ao menos uma prova escrita e trabalhos realizados e apresentados durante o semestre. O peso de cada atividade será definido segundo critérios do professor.Nota mínima de aprovação = 5,0" "A nota (N) será individual e será a média ponderada de componentes do projeto, tais como: Projeto Preliminar, Projeto Final, envolvimento do aluno com o projeto, Avaliação dos Pares, Apresentação de Trabalhos, dentre outros."

# Norma de recuperação
Replace-Text "Média Final = (N + Prova Recuperação)/2Nota Final mínima para aprovação= 5,0" "Média Final = (N + Prova Recuperação)/2"
